$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action")

# Switch active sheet from "Proxy" to "Action" (updates workbookView activeTab
# and moves tabSelected from the Proxy sheetView to the Action sheetView).
$ws.Select()

# Insert a new row 8 containing the "Action" label (moved up from F9), and
# shift the existing "Login" label from F10 up into F9, removing F10.
$ws.Range("F9").Copy($ws.Range("F8"))
$ws.Range("F10").Copy($ws.Range("F9"))
$ws.Range("F10").Clear()

# Update the selected cell on the Action sheet to I9.
$ws.Range("I9").Select()
